$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("A4").Value = "Research about Gannt Chart and How to Make them"
$ws.Range("B4").Value = 43700
$ws.Range("B4").NumberFormat = "d-mmm"
$ws.Range("C4").Value = 43701
$ws.Range("C4").NumberFormat = "d-mmm"
$ws.Range("D4").Value = 0.56944444444444442
$ws.Range("D4").NumberFormat = "h:mm AM/PM"
$ws.Range("E4").Value = 0.625
$ws.Range("E4").NumberFormat = "h:mm AM/PM"

# Row 5
$ws.Range("A5").Value = "Created Gantt Chart on Project Libre"
$ws.Range("B5").Value = 43701
$ws.Range("B5").NumberFormat = "d-mmm"
$ws.Range("C5").Value = 43702
$ws.Range("C5").NumberFormat = "d-mmm"
$ws.Range("D5").Value = 0.54166666666666663
$ws.Range("D5").NumberFormat = "h:mm AM/PM"
$ws.Range("E5").Value = 0.75
$ws.Range("E5").NumberFormat = "h:mm AM/PM"

# Move the active selection to A5, matching the saved view state
[void]$ws.Range("A5").Select()
